$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: Inflammatory-Mac)
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.065175
$ws.Range("N2").Value = 0.195525
$ws.Range("O2").Value = 0.009404016458916581
$ws.Range("P2").Value = 0.009404016458916581
$ws.Range("Q2").Value = 0.0594381227
$ws.Range("R2").Value = 0.5349431043
$ws.Range("S2").Value = 0.009404016458916581
$ws.Range("T2").Value = 0.009404016458916581

# Row 3 (Target cluster: Neutrophils)
$ws.Range("M3").Value = 6.718514333333332
$ws.Range("N3").Value = 20.155543
$ws.Range("O3").Value = 0.969405744075698
$ws.Range("P3").Value = 0.969405744075698
$ws.Range("Q3").Value = 6.12713278567511
$ws.Range("R3").Value = 55.144195071076
$ws.Range("S3").Value = 0.969405744075698
$ws.Range("T3").Value = 0.969405744075698

# Row 4 (Target cluster: Resolving-Mac)
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.14686
$ws.Range("N4").Value = 0.44058
$ws.Range("O4").Value = 0.02119023946538534
$ws.Range("P4").Value = 0.02119023946538533
$ws.Range("Q4").Value = 0.1339329911733333
$ws.Range("R4").Value = 1.20539692056
$ws.Range("S4").Value = 0.02119023946538534
$ws.Range("T4").Value = 0.02119023946538533
